# "Challenge Groups" workbook update
# - Fills in the previously-blank "country" cells for a handful of
#   participants (Group A2/A3, B2, C1, D1/D2/D3).
# - Moves the active selection to O18 (last cell the author was looking at).
#
# Cell values are written in the same left-to-right / top-to-bottom order
# the author typed them in, row by row, so that any newly-introduced
# shared-string entries land in the same order as the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group A2 / A3 (row 7) ---------------------------------------------
$ws.Range("O7").Value = "CIV"

# --- Group C1 (row 19) --------------------------------------------------
$ws.Range("E19").Value = "ETHIOPIA"

# --- Group A2 (row 7) ---------------------------------------------------
$ws.Range("J7").Value = "LIBANON"

# --- Group B2 (row 13) ---------------------------------------------------
$ws.Range("J13").Value = "TURKEY"

# --- Group D1 (row 25) ----------------------------------------------------
$ws.Range("E25").Value = "TUNISIA"

# --- Group D2 (row 25) ----------------------------------------------------
$ws.Range("J25").Value = "TUNISIA/IT"

# --- Group D3 (row 25) ----------------------------------------------------
$ws.Range("O25").Value = "MOROCO"

# --- Selection left where the author ended up editing --------------------
$ws.Range("O18").Select()
